$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Docentes responsáveis" (responsible faculty) entry needs to be added
# above the existing one (row 13), pushing that row and everything below it
# (old rows 13-22) down by one (new rows 14-23).
$ws.Rows(13).Insert()

# The row insert leaves a stray formatted-but-empty cell in column A (style
# carried over from the column); clear it since the new row only has B/C.
$ws.Range("A13").Clear()

# Copy the formatting (styles) from the existing docent row (now row 14,
# after the insert shifted it down) into the new row so B13/C13 match the
# look of B14/C14 exactly.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Fill in the new docent's name in both B and C, mirroring how the other
# docent row duplicates its text across columns B and C.
$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"
